$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gained two new quarterly data columns. Insert two blank
# columns before column D; the previously existing D:K data (8 quarters)
# shifts right to F:M, leaving D:E free for the two new quarters.
$ws.Range("D:E").Insert(-4161)

# Copy the number formatting from the (now shifted) former column D -
# i.e. current column F - onto the two new columns D:E, one contiguous
# data block at a time so blank separator/header rows (5,6,36,37,78,79)
# are left untouched.
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)

$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)

$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)

# Populate the two new columns with the new quarterly figures.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 35200
$ws.Range("E8").Value = 24500
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = -300
$ws.Range("E15").Value = -100
$ws.Range("D17").Value = 10100
$ws.Range("E17").Value = 7500
$ws.Range("D18").Value = 25100
$ws.Range("E18").Value = 17000
$ws.Range("D20").Value = -11400
$ws.Range("E20").Value = -6600
$ws.Range("D21").Value = 14400
$ws.Range("E21").Value = 10900
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 13700
$ws.Range("E23").Value = 10400
$ws.Range("D24").Value = 4700
$ws.Range("E24").Value = 2000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 9000
$ws.Range("E26").Value = 8300
$ws.Range("D27").Value = 9000
$ws.Range("E27").Value = 8300
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 500
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 11400
$ws.Range("E32").Value = 6600
$ws.Range("D33").Value = 9500
$ws.Range("E33").Value = 8300
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 9500
$ws.Range("E35").Value = 8300
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 147700
$ws.Range("E41").Value = 171600
$ws.Range("D42").Value = 10300
$ws.Range("E42").Value = 8300
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 17300
$ws.Range("E48").Value = 8100
$ws.Range("D49").Value = 66000
$ws.Range("E49").Value = 31100
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 4600
$ws.Range("E52").Value = 7300
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 2974000
$ws.Range("E54").Value = 2136800
$ws.Range("D57").Value = 21900
$ws.Range("E57").Value = 13200
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = "NA"
$ws.Range("E59").Value = "NA"
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 113200
$ws.Range("E61").Value = 53100
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 2599500
$ws.Range("E66").Value = 1841800
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 81600
$ws.Range("E72").Value = 73600
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 374500
$ws.Range("E76").Value = 295000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 9500
$ws.Range("E81").Value = 8300
$ws.Range("D83").Value = 800
$ws.Range("E83").Value = 500
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 59300
$ws.Range("E89").Value = -80400
$ws.Range("D91").Value = -400
$ws.Range("E91").Value = -1000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -71300
$ws.Range("E94").Value = -131800
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -11800
$ws.Range("E100").Value = 311000
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -23900
$ws.Range("E102").Value = 98800
